$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 329, shifting existing rows 329:435 down to 330:436
$ws.Rows.Item(329).Insert()

# Populate the newly inserted row 329 with the new record's data
$ws.Cells.Item(329, 1).Value = 10
$ws.Cells.Item(329, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(329, 3).Value = "La Araucanía"
$ws.Cells.Item(329, 4).Value = 44988
$ws.Cells.Item(329, 5).Value = 9
$ws.Cells.Item(329, 6).Value = 100114013
$ws.Cells.Item(329, 7).Value = "Zanahoria"
$ws.Cells.Item(329, 8).Value = "Sin especificar"
$ws.Cells.Item(329, 9).Value = "Primera"
$ws.Cells.Item(329, 10).Value = 100
$ws.Cells.Item(329, 11).Value = 7000
$ws.Cells.Item(329, 12).Value = 7000
$ws.Cells.Item(329, 13).Value = 7000
$ws.Cells.Item(329, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(329, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(329, 16).Value = 280
$ws.Cells.Item(329, 17).Value = 25
$ws.Cells.Item(329, 18).Value = "Hortaliza"
